$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellLines($row, $col, $lines) {
    $cell = $t.Cell($row, $col)
    $text = [string]::Join([string][char]11, $lines)
    $cell.Range.Text = $text
}

# Row 1
Set-CellLines 1 1 @("27 x 99", "  9    9", "  ----", "2|    |", "7|    |")
Set-CellLines 1 2 @("66 x 37", "  3    7", "  ----", "6|    |", "6|    |")
Set-CellLines 1 3 @("23 x 31", "  3    1", "  ----", "2|    |", "3|    |")

# Row 2
Set-CellLines 2 1 @("67 x 54", "  5    4", "  ----", "6|    |", "7|    |")
Set-CellLines 2 2 @("78 x 60", "  6    0", "  ----", "7|    |", "8|    |")
Set-CellLines 2 3 @("68 x 34", "  3    4", "  ----", "6|    |", "8|    |")

# Row 3
Set-CellLines 3 1 @("48 x 62", "  6    2", "  ----", "4|    |", "8|    |")
Set-CellLines 3 2 @("12 x 60", "  6    0", "  ----", "1|    |", "2|    |")
Set-CellLines 3 3 @("86 x 58", "  5    8", "  ----", "8|    |", "6|    |")

# Row 4
Set-CellLines 4 1 @("93 x 82", "  8    2", "  ----", "9|    |", "3|    |")
Set-CellLines 4 2 @("98 x 68", "  6    8", "  ----", "9|    |", "8|    |")
Set-CellLines 4 3 @("72 x 15", "  1    5", "  ----", "7|    |", "2|    |")

# Row 5
Set-CellLines 5 1 @("55 x 52", "  5    2", "  ----", "5|    |", "5|    |")
Set-CellLines 5 2 @("44 x 26", "  2    6", "  ----", "4|    |", "4|    |")
Set-CellLines 5 3 @("69 x 54", "  5    4", "  ----", "6|    |", "9|    |")
